$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Albahaca at Terminal La
# Palmera de La Serena. Insert a new row at the top of this data block
# (row 180), pushing the existing rows (180-200) down to (181-201), and
# populate the new row with the latest reading.
$ws.Rows(180).Insert()

$ws.Range("A180").Value = 8
$ws.Range("B180").Value = "Terminal La Palmera de La Serena"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 45124
$ws.Range("E180").Value = 4
$ws.Range("F180").Value = 100112052
$ws.Range("G180").Value = "Albahaca"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 1000
$ws.Range("K180").Value = 2800
$ws.Range("L180").Value = 3000
$ws.Range("M180").Value = 2900
$ws.Range("N180").Value = '$/paquete'
$ws.Range("O180").Value = "Región de Arica y Parinacota"
$ws.Range("P180").Value = 2900
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
